# Improve cell and inclusion segmentation algorithm:
# Update Number_of_Inclusions (column B) values and the recomputed
# Number_of_Inclusions_per_Nucleus (column D = B / C) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  B = 30  },
    @{ Row = 3;  B = 18  },
    @{ Row = 4;  B = 17  },
    @{ Row = 6;  B = 70  },
    @{ Row = 7;  B = 12  },
    @{ Row = 8;  B = 20  },
    @{ Row = 9;  B = 0   },
    @{ Row = 10; B = 0   },
    @{ Row = 12; B = 118 },
    @{ Row = 13; B = 26  },
    @{ Row = 14; B = 0   },
    @{ Row = 15; B = 0   },
    @{ Row = 19; B = 1   },
    @{ Row = 22; B = 3   },
    @{ Row = 23; B = 18  },
    @{ Row = 27; B = 0   },
    @{ Row = 29; B = 0   }
)

foreach ($u in $updates) {
    $r = $u.Row
    $b = $u.B
    $c = $ws.Cells.Item($r, 3).Value2

    $ws.Cells.Item($r, 2).Value = $b

    if ($c -ne 0) {
        $ws.Cells.Item($r, 4).Value = $b / $c
    } else {
        $ws.Cells.Item($r, 4).Value = 0
    }
}
